$d = $word.ActiveDocument

# The document body's 2nd paragraph holds a single run containing a VML
# <w:pict> with a <v:group> (a "canvas") that already has a couple of
# child shapes (two yellow rects, a connector and a textbox rect with
# "Shape text"). We need to add one more child shape to that group: a
# new textbox-style shape (v:shapetype #_x0000_t202 + v:shape) holding
# the text "Second shape", inserted right before the closing
# <w10:wrap type="none"/> of the group.
#
# VML shapes living inside a w:pict run aren't exposed through the
# Shapes/InlineShapes collections or through Find (they are not part of
# Range.Text), so the only reliable way to edit them is to replace the
# whole paragraph's XML via Range.InsertXML with an updated copy that
# contains the additional shapetype/shape.

$newParagraphXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:w10="urn:schemas-microsoft-com:office:word" w:rsidR="00DC1D4E" w:rsidRDefault="00566CC1" w:rsidP="00DC1D4E"><w:pPr><w:pStyle w:val="BodyTextIndent"/></w:pPr><w:r><w:pict><v:group id="_x0000_s1026" editas="canvas" style="width:342pt;height:180.65pt;mso-position-horizontal-relative:char;mso-position-vertical-relative:line" coordorigin="2785,-605" coordsize="6514,3468"><v:rect id="_x0000_s1028" style="position:absolute;left:2785;top:-260;width:2228;height:3123" fillcolor="yellow"/><v:rect id="_x0000_s1036" style="position:absolute;left:6765;top:-260;width:2200;height:2940" fillcolor="yellow"/><v:shapetype id="_x0000_t34" coordsize="21600,21600" o:spt="34" o:oned="t" adj="10800" path="m,l@0,0@0,21600,21600,21600e" filled="f"><v:stroke joinstyle="miter"/><v:formulas><v:f eqn="val #0"/></v:formulas><v:path arrowok="t" fillok="f" o:connecttype="none"/><v:handles><v:h position="#0,center"/></v:handles><o:lock v:ext="edit" shapetype="t"/></v:shapetype><v:shape id="_x0000_s1047" type="#_x0000_t34" style="position:absolute;left:4953;top:424;width:1870;height:917;rotation:180;flip:y"><v:stroke startarrow="block" endarrow="block"/></v:shape><v:rect id="_x0000_s1049" style="position:absolute;left:2859;top:1161;width:2094;height:358"><v:textbox><w:txbxContent><w:p w:rsidR="00DC1D4E" w:rsidRPr="00A524AD" w:rsidRDefault="00DC1D4E" w:rsidP="00DC1D4E"><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>Shape text</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect><v:shapetype id="_x0000_t202" coordsize="21600,21600" o:spt="202" path="m,l,21600r21600,l21600,xe"><v:stroke joinstyle="miter"/><v:path gradientshapeok="t" o:connecttype="rect"/></v:shapetype><v:shape id="_x0000_s1034" type="#_x0000_t202" style="position:absolute;left:2956;top:0;width:2057;height:345" stroked="f"><v:textbox><w:txbxContent><w:p><w:r><w:t>Second shape</w:t></w:r></w:p></w:txbxContent></v:textbox></v:shape><w10:wrap type="none"/></v:group></w:pict></w:r></w:p>
'@

$target = $d.Paragraphs(2).Range
[void]$target.InsertXML($newParagraphXml)
